$wb = $excel.ActiveWorkbook

# --- 1. Metadata sheet: bump the "Date" value ---------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2025-04-25T13:21:13+00:00"

# --- 2. "Mapping Table 4": insert a new mapping row for "Fourniture" ----
# Before:
#   row 8  -> Id_prescripteur            | related-to | MedicationRequest.requester
#   row 9  -> Identification_prescripteur| related-to | MedicationRequest.requester
# After:
#   row 8  -> Fourniture                 | equivalent | MedicationRequest.medication[x].extension.valueBoolean   (NEW)
#   row 9  -> Id_prescripteur            | related-to | MedicationRequest.requester
#   row 10 -> Identification_prescripteur| related-to | MedicationRequest.requester
$ws = $wb.Worksheets.Item("Mapping Table 4")

# Make room for a new row 10 by copying the formatting of row 9 down to it
# (row 10 does not exist yet, so it has no style of its own).
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the old rows 8 and 9 down into 9 and 10 (bottom-up so we do not
# clobber data before it has been copied).
$ws.Range("A10").Value2 = $ws.Range("A9").Value2
$ws.Range("C10").Value2 = $ws.Range("C9").Value2
$ws.Range("D10").Value2 = $ws.Range("D9").Value2

$ws.Range("A9").Value2 = $ws.Range("A8").Value2
$ws.Range("C9").Value2 = $ws.Range("C8").Value2
$ws.Range("D9").Value2 = $ws.Range("D8").Value2

# Write the new "Fourniture" mapping row into row 8.
$ws.Range("A8").Value2 = "Messages/M_prescription_médicaments/Prescription/Elément_prescr_médic/Fourniture"
$ws.Range("C8").Value2 = "equivalent"
$ws.Range("D8").Value2 = "MedicationRequest.medication[x].extension.valueBoolean"
